$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New registration row for "Julia"
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = "Julia"
$ws.Range("C2").Value = "j"

# Phone Number must be stored as text "12" (not a number). Build it via a
# formula that yields text, then flatten to a static value with a
# values-only paste so no extra number-format/style gets attached to the
# cell.
$ws.Range("D2").Formula = '=TEXT(12,"0")'
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)

$ws.Range("E2").Value = "j"
$ws.Range("F2").Value = "j"
